$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("r AnalysisUnit_Variable")

# Rename "CUSTOMER_BE_*" identifiers to "COUNTERPARTY_BIB_*" in columns B and C
# (rows 3-20) of the "r AnalysisUnit_Variable" sheet.
$ws2.Range("B3").Value = "COUNTERPARTY_BIB_SNDG"
$ws2.Range("C3").Value = "COUNTERPARTY_BIB_SNDG"

$ws2.Range("B4").Value = "COUNTERPARTY_BIB_TARGET"
$ws2.Range("C4").Value = "COUNTERPARTY_BIB_TARGET"

$ws2.Range("B5").Value = "COUNTERPARTY_BIB_IND_1"
$ws2.Range("C5").Value = "COUNTERPARTY_BIB_IND_1"

$ws2.Range("B6").Value = "COUNTERPARTY_BIB_IND_2"
$ws2.Range("C6").Value = "COUNTERPARTY_BIB_IND_2"

$ws2.Range("B7").Value = "COUNTERPARTY_BIB_IND_3"
$ws2.Range("C7").Value = "COUNTERPARTY_BIB_IND_3"

$ws2.Range("B8").Value = "COUNTERPARTY_BIB_IND_7"
$ws2.Range("C8").Value = "COUNTERPARTY_BIB_IND_7"

$ws2.Range("B9").Value = "COUNTERPARTY_BIB_IND_8"
$ws2.Range("C9").Value = "COUNTERPARTY_BIB_IND_8"

$ws2.Range("B10").Value = "COUNTERPARTY_BIB_IND_9"
$ws2.Range("C10").Value = "COUNTERPARTY_BIB_IND_9"

$ws2.Range("B11").Value = "COUNTERPARTY_BIB_IND_14"
$ws2.Range("C11").Value = "COUNTERPARTY_BIB_IND_14"

$ws2.Range("B12").Value = "COUNTERPARTY_BIB_IND_16"
$ws2.Range("C12").Value = "COUNTERPARTY_BIB_IND_16"

$ws2.Range("B13").Value = "COUNTERPARTY_BIB_IND_34"
$ws2.Range("C13").Value = "COUNTERPARTY_BIB_IND_34"

$ws2.Range("B14").Value = "COUNTERPARTY_BIB_IND_35"
$ws2.Range("C14").Value = "COUNTERPARTY_BIB_IND_35"

$ws2.Range("B15").Value = "COUNTERPARTY_BIB_IND_40"
$ws2.Range("C15").Value = "COUNTERPARTY_BIB_IND_40"

$ws2.Range("B16").Value = "COUNTERPARTY_BIB_IND_44"
$ws2.Range("C16").Value = "COUNTERPARTY_BIB_IND_44"

$ws2.Range("B17").Value = "COUNTERPARTY_BIB_IND_48"
$ws2.Range("C17").Value = "COUNTERPARTY_BIB_IND_48"

$ws2.Range("B18").Value = "COUNTERPARTY_BIB_IND_51"
$ws2.Range("C18").Value = "COUNTERPARTY_BIB_IND_51"

$ws2.Range("B19").Value = "COUNTERPARTY_BIB_IND_55"
$ws2.Range("C19").Value = "COUNTERPARTY_BIB_IND_55"

$ws2.Range("B20").Value = "COUNTERPARTY_BIB_IND_56"
$ws2.Range("C20").Value = "COUNTERPARTY_BIB_IND_56"

# Update the active sheet/selection: the "r AnalysisUnit_Variable" sheet
# becomes the active tab, with E3 selected (previously Analysis_Unit was
# the active tab).
$ws2.Activate()
$ws2.Range("E3").Select()
